$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 11, shifting existing rows 11:24 down to 12:25.
$ws.Rows.Item(11).Insert()

# Populate the newly-inserted row 11 with the new week's record.
$ws.Cells.Item(11, 1).Value = 10
$ws.Cells.Item(11, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(11, 3).Value = "La Araucanía"
$ws.Cells.Item(11, 4).Value = 44452
$ws.Cells.Item(11, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(11, 5).Value = 9
$ws.Cells.Item(11, 6).Value = 100114002
$ws.Cells.Item(11, 7).Value = "Camote"
$ws.Cells.Item(11, 8).Value = "Sin especificar"
$ws.Cells.Item(11, 9).Value = "Primera"
$ws.Cells.Item(11, 10).Value = 50
$ws.Cells.Item(11, 11).Value = 20000
$ws.Cells.Item(11, 12).Value = 20000
$ws.Cells.Item(11, 13).Value = 20000
$ws.Cells.Item(11, 14).Value = "$/malla 20 kilos"
$ws.Cells.Item(11, 15).Value = "Perú"
$ws.Cells.Item(11, 16).Value = 1000
$ws.Cells.Item(11, 17).Value = 20
$ws.Cells.Item(11, 18).Value = "Hortaliza"
